# Remove the "Appendix A: Latin Character Page" section (heading
# paragraph + the following include-directive paragraph) that was
# appended at the very end of the document, right after the
# "Acknowledgements" heading and its trailing blank paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that begins the Appendix A heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Appendix A: Latin Character Page")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Delete everything from the start of that heading paragraph through
    # the end of the document content (i.e. the heading paragraph and the
    # "<include=..\Examples\CharacterPageLatin.md>" paragraph that follows
    # it), leaving the preceding blank paragraph intact.
    $r = $d.Range($target.Range.Start, $d.Content.End)
    $r.Delete()
}
